$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 606.9091
$ws.Range("I15").Value = 606.9091
$ws.Range("K15").Value = 1820.7273
$ws.Range("M15").Value = -1651.7273
$ws.Range("H132").Value = 6474.48
$ws.Range("I132").Value = 6629.6313
$ws.Range("K132").Value = 19888.8939
$ws.Range("M132").Value = -17358.8939
$ws.Range("H135").Value = 1175.5294
$ws.Range("I135").Value = 732.26666
$ws.Range("K135").Value = 6590.39994
$ws.Range("M135").Value = -4055.39994
$ws.Range("H138").Value = 4367.9624
$ws.Range("J138").Value = 4920.512
$ws.Range("L138").Value = 14761.536
$ws.Range("N138").Value = -25041.536

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5210.923
$ws.Range("I2").Value = 4574.25
$ws.Range("K2").Value = 4574.25
$ws.Range("M2").Value = -4461.25
$ws.Range("H116").Value = 5210.923
$ws.Range("I116").Value = 4574.25
$ws.Range("K116").Value = 4574.25
$ws.Range("M116").Value = -2280.25
$ws.Range("H122").Value = 3258.7026
$ws.Range("I122").Value = 2407.3462
$ws.Range("K122").Value = 7222.0386
$ws.Range("M122").Value = -4772.0386
$ws.Range("H132").Value = 11042.333
$ws.Range("I132").Value = 2315.5
$ws.Range("K132").Value = 6946.5
$ws.Range("M132").Value = -4416.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5210.923
$ws.Range("I3").Value = 4574.25
$ws.Range("K3").Value = 4574.25
$ws.Range("M3").Value = -4460.25
$ws.Range("H20").Value = 2102.2632
$ws.Range("I20").Value = 1743.1818
$ws.Range("J20").Value = 2596
$ws.Range("K20").Value = 1743.1818
$ws.Range("L20").Value = 2596
$ws.Range("M20").Value = -1496.1818
$ws.Range("N20").Value = -3090
$ws.Range("H105").Value = 1546.0714
$ws.Range("I105").Value = 844
$ws.Range("J105").Value = 2248.1428
$ws.Range("K105").Value = 844
$ws.Range("L105").Value = 2248.1428
$ws.Range("M105").Value = 903
$ws.Range("N105").Value = -5742.1428
$ws.Range("H107").Value = 3167
$ws.Range("I107").Value = 2888.5715
$ws.Range("J107").Value = 3816.6667
$ws.Range("K107").Value = 2888.5715
$ws.Range("L107").Value = 3816.6667
$ws.Range("M107").Value = -968.5715
$ws.Range("N107").Value = -7656.6667
$ws.Range("H122").Value = 82666.664
$ws.Range("J122").Value = 82666.664
$ws.Range("L122").Value = 82666.664
$ws.Range("N122").Value = -92466.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 1812.125
$ws.Range("J15").Value = 2401.8
$ws.Range("L15").Value = 2401.8
$ws.Range("N15").Value = -2741.8
$ws.Range("H23").Value = 29000
$ws.Range("J23").Value = 28000
$ws.Range("L23").Value = 28000
$ws.Range("N23").Value = -28480
$ws.Range("H27").Value = 29000
$ws.Range("J27").Value = 28000
$ws.Range("L27").Value = 28000
$ws.Range("N27").Value = -28384
$ws.Range("H31").Value = 3407.9333
$ws.Range("J31").Value = 6661.6
$ws.Range("L31").Value = 6661.6
$ws.Range("N31").Value = -7251.6
$ws.Range("H34").Value = 3407.9333
$ws.Range("J34").Value = 6661.6
$ws.Range("L34").Value = 6661.6
$ws.Range("N34").Value = -7065.6
$ws.Range("H58").Value = 325853.72
$ws.Range("J58").Value = 6388.6665
$ws.Range("L58").Value = 6388.6665
$ws.Range("N58").Value = -6794.6665
$ws.Range("H129").Value = 63408.75
$ws.Range("J129").Value = 63408.75
$ws.Range("L129").Value = 63408.75
$ws.Range("N129").Value = -73408.75
$ws.Range("H136").Value = 325853.72
$ws.Range("J136").Value = 6388.6665
$ws.Range("L136").Value = 19165.9995
$ws.Range("N136").Value = -24265.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H43").Value = 300
$ws.Range("I43").Value = 300
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 900
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -786
$ws.Range("N43").ClearContents()
$ws.Range("H92").Value = 1542.5714
$ws.Range("I92").Value = 716.3333
$ws.Range("J92").Value = 2162.25
$ws.Range("K92").Value = 2148.9999
$ws.Range("L92").Value = 6486.75
$ws.Range("M92").Value = -900.9998999999998
$ws.Range("N92").Value = -8982.75
$ws.Range("H99").Value = 4079.5
$ws.Range("I99").Value = 2300
$ws.Range("J99").Value = 4672.6665
$ws.Range("K99").Value = 6900
$ws.Range("L99").Value = 14017.9995
$ws.Range("M99").Value = -4654
$ws.Range("N99").Value = -18509.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1908.4445
$ws.Range("I13").Value = 2290
$ws.Range("K13").Value = 2290
$ws.Range("M13").Value = -2151
$ws.Range("H40").Value = 11499.5
$ws.Range("J40").Value = 11499.5
$ws.Range("L40").Value = 11499.5
$ws.Range("N40").Value = -11801.5
$ws.Range("H126").Value = 4657.75
$ws.Range("I126").Value = 2881
$ws.Range("K126").Value = 8643
$ws.Range("M126").Value = -6173
$ws.Range("H132").Value = 292312.78
$ws.Range("I132").Value = 359531.53
$ws.Range("J132").Value = 104100.3
$ws.Range("K132").Value = 1078594.59
$ws.Range("L132").Value = 312300.9
$ws.Range("M132").Value = -1076064.59
$ws.Range("N132").Value = -317360.9
$ws.Range("H133").Value = 67128.57000000001
$ws.Range("J133").Value = 67128.57000000001
$ws.Range("L133").Value = 67128.57000000001
$ws.Range("N133").Value = -77248.57000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 38465920
$ws.Range("I122").Value = 55558548
$ws.Range("K122").Value = 166675644
$ws.Range("M122").Value = -166673194
$ws.Range("H126").Value = 3591.68
$ws.Range("I126").Value = 2486.2666
$ws.Range("K126").Value = 7458.7998
$ws.Range("M126").Value = -4988.7998
$ws.Range("H128").Value = 60500
$ws.Range("I128").Value = 60000
$ws.Range("J128").Value = 61000
$ws.Range("K128").Value = 60000
$ws.Range("L128").Value = 61000
$ws.Range("M128").Value = -55020
$ws.Range("N128").Value = -70960
$ws.Range("H132").Value = 44639.44
$ws.Range("I132").Value = 2499.0715
$ws.Range("K132").Value = 7497.2145
$ws.Range("M132").Value = -4967.2145
$ws.Range("H133").Value = 64994.25
$ws.Range("J133").Value = 64994.25
$ws.Range("L133").Value = 64994.25
$ws.Range("N133").Value = -75114.25
$ws.Range("H136").Value = 378707.2
$ws.Range("I136").Value = 437419
$ws.Range("J136").Value = 228665.89
$ws.Range("K136").Value = 1312257
$ws.Range("L136").Value = 685997.67
$ws.Range("M136").Value = -1309707
$ws.Range("N136").Value = -691097.67
